$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 4 (shifts current rows 4-19 down to 6-21)
$ws.Range("A4:A5").EntireRow.Insert()

# Insert one more new row before what is now row 7 (the old row 5, "+5522981222545")
$ws.Range("A7").EntireRow.Insert()

# The source data stores phone numbers, DDD codes, and dates as plain text
# (not numbers / dates), so force text format on the new cells before
# writing their values to avoid Excel auto-converting them.
$newRows = @(4, 5, 7)
foreach ($r in $newRows) {
    $ws.Range("A" + $r + ":C" + $r).NumberFormat = "@"
}

# Fill in the three newly inserted rows with their values
$ws.Range("A4").Value = "+5519997201600"
$ws.Range("B4").Value = "19"
$ws.Range("C4").Value = "2024-10-18"

$ws.Range("A5").Value = "+5511975292030"
$ws.Range("B5").Value = "11"
$ws.Range("C5").Value = "2024-10-14"

$ws.Range("A7").Value = "+5515996313442"
$ws.Range("B7").Value = "15"
$ws.Range("C7").Value = "2024-10-09"
